$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 14.33087667279843

# Row 3
$ws.Range("B3").Value = 0.04172184405617529
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 18.71679738969934
$ws.Range("E3").Value = 2797.565817734744
$ws.Range("G3").Value = 2816.629228217133
